$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update marking per correct answer (row 11 "Marking", column B "Right")
$ws.Range("B11").Value = 5

# Update total correct marks (row 12 "Total", column B "Right")
$ws.Range("B12").Value = 120

# Update the correct/total marks display text (row 12, column E "Max")
$ws.Range("E12").Value = "120/140"
